# The edit inserts one new price-record row into the weekly "Jengibre" (ginger)
# price table for "Vega Modelo de Temuco". This shifts every existing data row
# from row 53 down through row 169 to rows 54 through 170, and populates the
# newly opened row 53 with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 53; Excel shifts rows 53:169 down to 54:170
# and extends the used range to row 170 automatically.
$ws.Rows.Item(53).Insert()

# Seed the new row 53 with the same layout/values as the row right below it
# (which is the row that used to be row 53), then overwrite just the cells
# that differ for the new record.
$ws.Range("A54:R54").Copy()
$ws.Range("A53:R53").PasteSpecial()

$ws.Range("D53").Value = 44720
$ws.Range("K53").Value = 20000
$ws.Range("L53").Value = 20000
$ws.Range("M53").Value = 20000
$ws.Range("P53").Value = 1538
